$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 583, pushing the existing rows 583-629
# (and their formatting) down to 584-630.
$ws.Rows(583).Insert()

# Populate the newly inserted row 583 with the new weekly price record.
$ws.Range("A583").Value = 3
$ws.Range("B583").Value = "Femacal de La Calera"
$ws.Range("C583").Value = "Coquimbo"
$ws.Range("D583").Value = 45223
$ws.Range("E583").Value = 5
$ws.Range("F583").Value = 100112012
$ws.Range("G583").Value = "Espinaca"
$ws.Range("H583").Value = "Sin especificar"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 170
$ws.Range("K583").Value = 4000
$ws.Range("L583").Value = 4300
$ws.Range("M583").Value = 4141
$ws.Range("N583").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O583").Value = "Provincia de Quillota"
$ws.Range("P583").Value = 1380
$ws.Range("Q583").Value = 3
$ws.Range("R583").Value = "Hortaliza"
